# Auto commit at 2025-12-04 7:37:48.87
# Updates the "Metrics" sheet's source figures (B2:B13) with refreshed
# totals, and leaves the downstream "today" sheet's formulas (which pull
# from Metrics!) to recalculate automatically. Also refreshes the saved
# cell selections on both sheets.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# --- Metrics!B2:B13 -> new values -----------------------------------
$metrics.Range("B2").Value  = 38769.18
$metrics.Range("B3").Value  = 33461.71
$metrics.Range("B4").Value  = 11840.63
$metrics.Range("B5").Value  = 1595
$metrics.Range("B6").Value  = 5241476.290000001
$metrics.Range("B7").Value  = 4433814.6700000009
$metrics.Range("B8").Value  = 1543797.5100000005
$metrics.Range("B9").Value  = 204302
$metrics.Range("B10").Value = 33706857.280000001
$metrics.Range("B11").Value = 31709089.830000002
$metrics.Range("B12").Value = 11825519.549999997
$metrics.Range("B13").Value = 1301932

# --- refresh the stored selections ------------------------------------
# Move to Metrics!F19 first (matches the saved selection for that sheet).
$metrics.Range("F19").Select()

# Re-activate "today" (it is the tab that was showing before this edit)
# and leave its selection on F5, matching the saved selection there.
$today.Activate()
$today.Range("F5").Select()
